# Loop through the people in the list, look at the questions each person got
# wrong (column D, "Q wrong", e.g. "3, 5"), pull the corresponding message
# text for each wrong question out of the "Messages for each question" block
# (column E, row = question number + 1, since row 2 holds the message for
# question 1, row 3 for question 2, etc.) and build the combined message for
# column F ("List of message").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new wrap-text style exists / gets applied to the output cells.
$firstDataRow = 2
$idCol = "A"              # ID numbers - one row per person
$wrongCol = "D"           # Q wrong
$outCol = "F"             # List of message

$row = $firstDataRow
while ($ws.Range($idCol + $row).Value2 -ne $null -and $ws.Range($idCol + $row).Value2 -ne "") {

    $wrongText = [string]$ws.Range($wrongCol + $row).Value2
    $questionNumbers = $wrongText -split ","

    $combined = ""
    foreach ($q in $questionNumbers) {
        $qNum = $q.Trim()
        if ($qNum -ne "") {
            $msgRow = [int]$qNum + 1
            $msgText = [string]$ws.Range("E" + $msgRow).Value2
            $combined = $combined + $msgText + "  " + [char]10
        }
    }

    $outCell = $ws.Range($outCol + $row)
    $outCell.Value = $combined
    $outCell.WrapText = $true

    $ws.Rows.Item($row).RowHeight = 21.6 * $questionNumbers.Count

    $row = $row + 1
}

$null = $ws.Range("F2").Select()
